$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.537.07'
$ws.Range('E2').Value = '  +0.70%  '
$ws.Range('D3').Value = '1.729.67'
$ws.Range('E3').Value = '  +0.66%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9995'
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '245.32'
$ws.Range('E5').Value = '  +2.82%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4806'
$ws.Range('E7').Value = '  +1.63%  '
$ws.Range('E8').Value = '  +1.40%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06222'
$ws.Range('E9').Value = '  +0.19%  '
$ws.Range('D10').Value = '1.727.46'
$ws.Range('E10').Value = '  +0.55%  '
$ws.Range('E11').Value = '  +1.36%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '15.72'
$ws.Range('E12').Value = '  +2.53%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.6164'
$ws.Range('E13').Value = '  +4.18%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.542'
$ws.Range('E14').Value = '  +2.94%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '77.21'
$ws.Range('E15').Value = '  +1.33%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.9999'
$ws.Range('E16').Value = '  -0.06%  '
$ws.Range('D17').Value = '26.535.13'
$ws.Range('E17').Value = '  +0.71%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '1.000'
$ws.Range('E18').Value = '  -0.04%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000006955'
$ws.Range('E19').Value = '  +2.11%  '
$ws.Range('D21').Value = '1.949.86'
$ws.Range('E21').Value = '  +0.66%  '
$ws.Range('E22').Value = '  -0.35%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '8.932'
$ws.Range('E23').Value = '  +1.95%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.289'
$ws.Range('E24').Value = '  -0.54%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '136.57'
$ws.Range('E25').Value = '  +1.19%  '
$ws.Range('E26').Value = '  +0.78%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.799'
$ws.Range('E27').Value = '  +2.60%  '
$ws.Range('E28').Value = '  +0.05%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '106.85'
$ws.Range('E29').Value = '  -1.22%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '3.993'
$ws.Range('E30').Value = '  -0.24%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.07990'
$ws.Range('E31').Value = '  +3.21%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.721'
$ws.Range('E32').Value = '  +0.96%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.04593'
$ws.Range('E33').Value = '  +3.61%  '
$ws.Range('B34').Value = 'HuobiToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.617'
$ws.Range('E34').Value = '  +0.10%  '
$ws.Range('B35').Value = 'ImmutableX'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.6368'
$ws.Range('E35').Value = '  +2.86%  '
$ws.Range('B36').Value = 'ARBITRUM'
$ws.Range('C36').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.9962'
$ws.Range('E36').Value = '  +1.71%  '
$ws.Range('B37').Value = 'TrustWalletToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.9211'
$ws.Range('E37').Value = '  -1.36%  '
$ws.Range('B38').Value = 'RenderToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.094'
$ws.Range('E38').Value = '  +9.44%  '
$ws.Range('B39').Value = 'MXToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.405'
$ws.Range('E39').Value = '  -0.41%  '
$ws.Range('B40').Value = 'Quant'
$ws.Range('C40').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '104.81'
$ws.Range('E40').Value = '  -7.95%  '
$ws.Range('B41').Value = 'PaxDollar'
$ws.Range('C41').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.003'
$ws.Range('E41').Value = '  +0.28%  '
$ws.Range('B42').Value = 'VeChain'
$ws.Range('C42').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.01508'
$ws.Range('E42').Value = '  +2.22%  '
$ws.Range('B43').Value = 'FraxShare'
$ws.Range('C43').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.586'
$ws.Range('E43').Value = '  +4.87%  '
$ws.Range('B44').Value = 'TheSandbox'
$ws.Range('C44').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.3899'
$ws.Range('E44').Value = '  +2.28%  '
$ws.Range('B45').Value = 'Aptos'
$ws.Range('C45').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '7.001'
$ws.Range('E45').Value = '  +11.30%  '
$ws.Range('B46').Value = 'Algorand'
$ws.Range('C46').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.1184'
$ws.Range('E46').Value = '  +1.32%  '
$ws.Range('B47').Value = 'Cronos'
$ws.Range('C47').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.05339'
$ws.Range('E47').Value = '  +1.03%  '
$ws.Range('B48').Value = 'Elrond'
$ws.Range('C48').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '30.98'
$ws.Range('E48').Value = '  +1.88%  '
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.877'
$ws.Range('E49').Value = '  +2.36%  '
$ws.Range('B50').Value = 'NEARProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.266'
$ws.Range('E50').Value = '  +4.09%  '
$ws.Range('B51').Value = 'Decentraland'
$ws.Range('C51').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.3430'
$ws.Range('E51').Value = '  +1.80%  '
